$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 16:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3222123
$ws.Range("C4").Value = 2124
$ws.Range("E4").Value = 1659641
$ws.Range("G4").Value = 47
$ws.Range("H4").Value = 135869

# India (row 6)
$ws.Range("B6").Value = 798161
$ws.Range("C6").Value = 3319
$ws.Range("D6").Value = 497690
$ws.Range("E6").Value = 278815

# Arabia Saudita (row 17)
$ws.Range("B17").Value = 226486
$ws.Range("C17").Value = 3159
$ws.Range("D17").Value = 163026
$ws.Range("E17").Value = 61309
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 2151

# Alemania (row 19)
$ws.Range("B19").Value = 199254
$ws.Range("C19").Value = 56
$ws.Range("E19").Value = 6528
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 9126

# Argentina (row 25)
$ws.Range("D25").Value = 38984
$ws.Range("E25").Value = 49960
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 1749

# Emiratos Arabes Unidos (row 37)
$ws.Range("B37").Value = 54050
$ws.Range("C37").Value = 473
$ws.Range("D37").Value = 43969
$ws.Range("E37").Value = 9751
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 330

# Paises Bajos (row 40)
$ws.Range("B40").Value = 50840
$ws.Range("C40").Value = 42
$ws.Range("H40").Value = 6136

# Azerbaiyan (row 58)
$ws.Range("B58").Value = 22990
$ws.Range("C58").Value = 526
$ws.Range("D58").Value = 14093
$ws.Range("E58").Value = 8605
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 292

# Kenia (row 76)
$ws.Range("B76").Value = 9448
$ws.Range("C76").Value = 473
$ws.Range("D76").Value = 2733
$ws.Range("E76").Value = 6534
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 8
$ws.Range("H76").Value = 181

# Republica de Macedonia (row 82)
$ws.Range("B82").Value = 7777
$ws.Range("C82").Value = 205
$ws.Range("D82").Value = 3960
$ws.Range("E82").Value = 3449
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 368

# Bosnia y Herzegovina (row 89)
$ws.Range("B89").Value = 6402
$ws.Range("C89").Value = 316
$ws.Range("D89").Value = 3037
$ws.Range("E89").Value = 3149
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 216

# Malaui (row 116)
$ws.Range("E116").Value = 1586
$ws.Range("G116").Value = 4
$ws.Range("H116").Value = 29

# Tunez (row 131)
$ws.Range("B131").Value = 1240
$ws.Range("C131").Value = 9
$ws.Range("D131").Value = 1067
$ws.Range("E131").Value = 123

# Liberia (row 144)
$ws.Range("B144").Value = 963
$ws.Range("C144").Value = 6
$ws.Range("D144").Value = 400
$ws.Range("E144").Value = 516
$ws.Range("G144").Value = 5
$ws.Range("H144").Value = 47

# Principado de Andorra (row 147)
$ws.Range("D147").Value = 803
$ws.Range("E147").Value = 0

# Birmania (row 164)
$ws.Range("B164").Value = 326
$ws.Range("C164").Value = 7
$ws.Range("D164").Value = 256
$ws.Range("E164").Value = 64

# Re-sort the data block (A3:H219, with row 3 as header) descending by
# column B ("Casos totales") to reflect each country's latest totals.
$dataRange = $ws.Range("A3:H219")
$keyCol = $ws.Range("B4")
$dataRange.Sort($keyCol, 2, $null, $null, 1, $null, 1, 1)
